$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 518.2
$ws.Cells.Item(28, 9).Value = 409.1111
$ws.Cells.Item(28, 10).Value = 1500
$ws.Cells.Item(28, 11).Value = 409.1111
$ws.Cells.Item(28, 12).Value = 1500
$ws.Cells.Item(28, 13).Value = 75.88889999999998
$ws.Cells.Item(28, 14).Value = -2470
$ws.Cells.Item(32, 8).Value = 3875.5293
$ws.Cells.Item(32, 10).Value = 4126.143
$ws.Cells.Item(32, 12).Value = 4126.143
$ws.Cells.Item(32, 14).Value = -4778.143
$ws.Cells.Item(40, 8).Value = 3631.257
$ws.Cells.Item(40, 9).Value = 3993.2173
$ws.Cells.Item(40, 11).Value = 3993.2173
$ws.Cells.Item(40, 13).Value = -3818.2173
$ws.Cells.Item(98, 8).Value = 1459.7073
$ws.Cells.Item(98, 9).Value = 1503.3158
$ws.Cells.Item(98, 11).Value = 1503.3158
$ws.Cells.Item(98, 13).Value = -5.315800000000081
$ws.Cells.Item(103, 8).Value = 1304
$ws.Cells.Item(103, 10).Value = 1300
$ws.Cells.Item(103, 12).Value = 3900
$ws.Cells.Item(103, 14).Value = -5072
$ws.Cells.Item(122, 8).Value = 1459.7073
$ws.Cells.Item(122, 9).Value = 1503.3158
$ws.Cells.Item(122, 11).Value = 4509.9474
$ws.Cells.Item(122, 13).Value = -2059.9474
$ws.Cells.Item(129, 8).Value = 2596.7896
$ws.Cells.Item(129, 9).Value = 2401.4666
$ws.Cells.Item(129, 11).Value = 7204.399800000001
$ws.Cells.Item(129, 13).Value = -2204.399800000001
$ws.Cells.Item(131, 8).Value = 489.44446
$ws.Cells.Item(131, 9).Value = 489.44446
$ws.Cells.Item(131, 11).Value = 1468.33338
$ws.Cells.Item(131, 13).Value = 3571.66662
$ws.Cells.Item(132, 8).Value = 322208.34
$ws.Cells.Item(132, 9).Value = 340055.1
$ws.Cells.Item(132, 11).Value = 1020165.3
$ws.Cells.Item(132, 13).Value = -1017635.3
$ws.Cells.Item(138, 8).Value = 3640.6
$ws.Cells.Item(138, 10).Value = 3816.3
$ws.Cells.Item(138, 12).Value = 11448.9
$ws.Cells.Item(138, 14).Value = -21728.9
$ws.Cells.Item(141, 8).Value = 1266.4546
$ws.Cells.Item(141, 9).Value = 1193.1
$ws.Cells.Item(141, 11).Value = 3579.3
$ws.Cells.Item(141, 13).Value = 1600.7
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 1356
$ws.Cells.Item(25, 9).Value = 1070
$ws.Cells.Item(25, 11).Value = 1070
$ws.Cells.Item(25, 13).Value = -668
$ws.Cells.Item(32, 8).Value = 5161810.5
$ws.Cells.Item(32, 9).Value = 6497523.5
$ws.Cells.Item(32, 10).Value = 19316.85
$ws.Cells.Item(32, 11).Value = 6497523.5
$ws.Cells.Item(32, 12).Value = 19316.85
$ws.Cells.Item(32, 13).Value = -6497236.5
$ws.Cells.Item(32, 14).Value = -19890.85
$ws.Cells.Item(43, 8).Value = 34152
$ws.Cells.Item(43, 9).Value = 34170.5
$ws.Cells.Item(43, 10).Value = 34141.43
$ws.Cells.Item(43, 11).Value = 34170.5
$ws.Cells.Item(43, 12).Value = 34141.43
$ws.Cells.Item(43, 13).Value = -33857.5
$ws.Cells.Item(43, 14).Value = -34767.43
$ws.Cells.Item(74, 8).Value = 4313763.5
$ws.Cells.Item(74, 9).Value = 5683582.5
$ws.Cells.Item(74, 10).Value = 8617.143
$ws.Cells.Item(74, 11).Value = 5683582.5
$ws.Cells.Item(74, 12).Value = 8617.143
$ws.Cells.Item(74, 13).Value = -5682708.5
$ws.Cells.Item(74, 14).Value = -10365.143
$ws.Cells.Item(77, 8).Value = 4313763.5
$ws.Cells.Item(77, 9).Value = 5683582.5
$ws.Cells.Item(77, 10).Value = 8617.143
$ws.Cells.Item(77, 11).Value = 28417912.5
$ws.Cells.Item(77, 12).Value = 43085.715
$ws.Cells.Item(77, 13).Value = -28413544.5
$ws.Cells.Item(77, 14).Value = -51821.715
$ws.Cells.Item(102, 8).Value = 6299.727
$ws.Cells.Item(102, 9).Value = 6588.6665
$ws.Cells.Item(102, 11).Value = 6588.6665
$ws.Cells.Item(102, 13).Value = -4966.6665
$ws.Cells.Item(122, 8).Value = 2633.4119
$ws.Cells.Item(122, 9).Value = 2197.7856
$ws.Cells.Item(122, 11).Value = 6593.3568
$ws.Cells.Item(122, 13).Value = -4143.3568
$ws.Cells.Item(132, 8).Value = 690156.8
$ws.Cells.Item(132, 10).Value = 11321.667
$ws.Cells.Item(132, 12).Value = 33965.001
$ws.Cells.Item(132, 14).Value = -39025.001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 900.76666
$ws.Cells.Item(94, 9).Value = 742.7222
$ws.Cells.Item(94, 10).Value = 1137.8334
$ws.Cells.Item(94, 11).Value = 742.7222
$ws.Cells.Item(94, 12).Value = 1137.8334
$ws.Cells.Item(94, 13).Value = -291.7222
$ws.Cells.Item(94, 14).Value = -2039.8334
$ws.Cells.Item(105, 8).Value = 4386.375
$ws.Cells.Item(105, 9).Value = 4049.8
$ws.Cells.Item(105, 10).Value = 4947.3335
$ws.Cells.Item(105, 11).Value = 4049.8
$ws.Cells.Item(105, 12).Value = 4947.3335
$ws.Cells.Item(105, 13).Value = -2302.8
$ws.Cells.Item(105, 14).Value = -8441.333500000001
$ws.Cells.Item(134, 8).Value = 657647.0600000001
$ws.Cells.Item(134, 9).Value = 889281.9399999999
$ws.Cells.Item(134, 11).Value = 2667845.82
$ws.Cells.Item(134, 13).Value = -2665310.82
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8760.5
$ws.Cells.Item(31, 9).Value = 1388.762
$ws.Cells.Item(31, 10).Value = 13451.606
$ws.Cells.Item(31, 11).Value = 1388.762
$ws.Cells.Item(31, 12).Value = 13451.606
$ws.Cells.Item(31, 13).Value = -1093.762
$ws.Cells.Item(31, 14).Value = -14041.606
$ws.Cells.Item(34, 8).Value = 8760.5
$ws.Cells.Item(34, 9).Value = 1388.762
$ws.Cells.Item(34, 10).Value = 13451.606
$ws.Cells.Item(34, 11).Value = 1388.762
$ws.Cells.Item(34, 12).Value = 13451.606
$ws.Cells.Item(34, 13).Value = -1186.762
$ws.Cells.Item(34, 14).Value = -13855.606
$ws.Cells.Item(50, 8).Value = 119995
$ws.Cells.Item(50, 10).Value = 119995
$ws.Cells.Item(50, 12).Value = 119995
$ws.Cells.Item(50, 14).Value = -121245
$ws.Cells.Item(60, 8).Value = 75355.22
$ws.Cells.Item(60, 10).Value = 79837.125
$ws.Cells.Item(60, 12).Value = 79837.125
$ws.Cells.Item(60, 14).Value = -80859.125
$ws.Cells.Item(62, 8).Value = 4547.5713
$ws.Cells.Item(62, 9).Value = 4298.8887
$ws.Cells.Item(62, 11).Value = 4298.8887
$ws.Cells.Item(62, 13).Value = -3674.8887
$ws.Cells.Item(65, 8).Value = 4547.5713
$ws.Cells.Item(65, 9).Value = 4298.8887
$ws.Cells.Item(65, 11).Value = 21494.4435
$ws.Cells.Item(65, 13).Value = -18374.4435
$ws.Cells.Item(105, 8).Value = 25242.934
$ws.Cells.Item(105, 9).Value = 26831.715
$ws.Cells.Item(105, 11).Value = 26831.715
$ws.Cells.Item(105, 13).Value = -25084.715
$ws.Cells.Item(132, 8).Value = 7826768
$ws.Cells.Item(132, 9).Value = 16526.875
$ws.Cells.Item(132, 11).Value = 49580.625
$ws.Cells.Item(132, 13).Value = -47050.625
$ws.Cells.Item(134, 8).Value = 6652.6924
$ws.Cells.Item(134, 9).Value = 2215.5908
$ws.Cells.Item(134, 11).Value = 6646.7724
$ws.Cells.Item(134, 13).Value = -4111.7724
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 73263950
$ws.Cells.Item(4, 9).Value = 84221040
$ws.Cells.Item(4, 11).Value = 252663120
$ws.Cells.Item(4, 13).Value = -252663008
$ws.Cells.Item(14, 8).Value = 118206.65
$ws.Cells.Item(14, 9).Value = 118206.65
$ws.Cells.Item(14, 11).Value = 354619.95
$ws.Cells.Item(14, 13).Value = -354446.95
$ws.Cells.Item(121, 8).Value = 917.5
$ws.Cells.Item(121, 9).Value = 734.3333
$ws.Cells.Item(121, 10).Value = 1027.4
$ws.Cells.Item(121, 11).Value = 2202.9999
$ws.Cells.Item(121, 12).Value = 3082.2
$ws.Cells.Item(121, 13).Value = -892.9998999999998
$ws.Cells.Item(121, 14).Value = -5702.200000000001
$ws.Cells.Item(122, 8).Value = 771.1316
$ws.Cells.Item(122, 9).Value = 495.7143
$ws.Cells.Item(122, 10).Value = 833.3226
$ws.Cells.Item(122, 11).Value = 4461.428699999999
$ws.Cells.Item(122, 12).Value = 7499.903399999999
$ws.Cells.Item(122, 13).Value = -2011.428699999999
$ws.Cells.Item(122, 14).Value = -12399.9034
$ws.Cells.Item(131, 8).Value = 7274.425
$ws.Cells.Item(131, 9).Value = 992.5
$ws.Cells.Item(131, 10).Value = 8383
$ws.Cells.Item(131, 11).Value = 2977.5
$ws.Cells.Item(131, 12).Value = 25149
$ws.Cells.Item(131, 13).Value = 2062.5
$ws.Cells.Item(131, 14).Value = -35229
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2976.5
$ws.Cells.Item(102, 9).Value = 2418.611
$ws.Cells.Item(102, 10).Value = 7997.5
$ws.Cells.Item(102, 11).Value = 2418.611
$ws.Cells.Item(102, 12).Value = 7997.5
$ws.Cells.Item(102, 13).Value = -796.6109999999999
$ws.Cells.Item(102, 14).Value = -11241.5
$ws.Cells.Item(113, 8).Value = 1265.6471
$ws.Cells.Item(113, 9).Value = 1251.0667
$ws.Cells.Item(113, 10).Value = 1375
$ws.Cells.Item(113, 11).Value = 1251.0667
$ws.Cells.Item(113, 12).Value = 1375
$ws.Cells.Item(113, 13).Value = 918.9332999999999
$ws.Cells.Item(113, 14).Value = -5715
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(11, 8).Value = 5263
$ws.Cells.Item(11, 10).Value = 5263
$ws.Cells.Item(11, 12).Value = 5263
$ws.Cells.Item(11, 14).Value = -5543
$ws.Cells.Item(46, 8).Value = 2675.3794
$ws.Cells.Item(46, 9).Value = 1240.5
$ws.Cells.Item(46, 10).Value = 3688.2354
$ws.Cells.Item(46, 11).Value = 1240.5
$ws.Cells.Item(46, 12).Value = 3688.2354
$ws.Cells.Item(46, 13).Value = -1052.5
$ws.Cells.Item(46, 14).Value = -4064.2354
$ws.Cells.Item(61, 8).Value = 11321.154
$ws.Cells.Item(61, 9).Value = 14573.368
$ws.Cells.Item(61, 10).Value = 2493.7144
$ws.Cells.Item(61, 11).Value = 14573.368
$ws.Cells.Item(61, 12).Value = 2493.7144
$ws.Cells.Item(61, 13).Value = -14371.368
$ws.Cells.Item(61, 14).Value = -2897.7144
$ws.Cells.Item(93, 8).Value = 5356.4
$ws.Cells.Item(93, 9).Value = 3034.3
$ws.Cells.Item(93, 11).Value = 3034.3
$ws.Cells.Item(93, 13).Value = -1786.3
$ws.Cells.Item(110, 8).Value = 83081.336
$ws.Cells.Item(110, 10).Value = 83081.336
$ws.Cells.Item(110, 12).Value = 83081.336
$ws.Cells.Item(110, 14).Value = -91261.336
$ws.Cells.Item(113, 8).Value = 11321.154
$ws.Cells.Item(113, 9).Value = 14573.368
$ws.Cells.Item(113, 10).Value = 2493.7144
$ws.Cells.Item(113, 11).Value = 14573.368
$ws.Cells.Item(113, 12).Value = 2493.7144
$ws.Cells.Item(113, 13).Value = -12403.368
$ws.Cells.Item(113, 14).Value = -6833.7144
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 13).Value = ""
$ws.Cells.Item(107, 8).Value = 4067.077
$ws.Cells.Item(107, 9).Value = 2235
$ws.Cells.Item(107, 11).Value = 6705
$ws.Cells.Item(107, 13).Value = -4785
$ws.Cells.Item(132, 8).Value = 20328612
$ws.Cells.Item(132, 9).Value = 1921863.6
$ws.Cells.Item(132, 11).Value = 5765590.800000001
$ws.Cells.Item(132, 13).Value = -5763060.800000001
$ws.Cells.Item(136, 8).Value = 16563778
$ws.Cells.Item(136, 9).Value = 19047096
$ws.Cells.Item(136, 11).Value = 57141288
$ws.Cells.Item(136, 13).Value = -57138738
